# Actualización automática 2025-05-30 08:25:07
#
# Inserts a new detail row (PAREDES AGUILAR ILIANA GABRIELA / CONSTRUCCION,
# INGENIERIA Y TECNOLOGIA CONSTRUINTEC SAS) at row 241 on both data sheets,
# pushing the existing rows 241-261 down to 242-262, and refreshes the
# trailing "summary" rows (the "x de N" counters on "VENTAS POR GRUPO" and
# the plain numeric totals on "VENTA MENSUAL") to their new row numbers.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "VENTAS POR GRUPO" (columns A:N) ----------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(241).Insert()

$ws1.Cells.Item(241, 1).Value = "PAREDES AGUILAR ILIANA GABRIELA"
$ws1.Cells.Item(241, 2).Value = "CONSTRUCCION, INGENIERIA Y TECNOLOGIA CONSTRUINTEC SAS"
for ($col = 3; $col -le 14; $col++) {
    $ws1.Cells.Item(241, $col).Value = 0
}

# the "x de 260" labels (row 262 before the insert) are now on row 263;
# bump the denominator to 261 while keeping the numerator untouched.
for ($col = 3; $col -le 14; $col++) {
    $cell = $ws1.Cells.Item(263, $col)
    [string]$text = $cell.Text
    $cell.Value = $text.Replace("de 260", "de 261")
}

# ---- Sheet 2: "VENTA MENSUAL" (columns A:F) --------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(241).Insert()

$ws2.Cells.Item(241, 1).Value = "PAREDES AGUILAR ILIANA GABRIELA"
$ws2.Cells.Item(241, 2).Value = "CONSTRUCCION, INGENIERIA Y TECNOLOGIA CONSTRUINTEC SAS"
for ($col = 3; $col -le 6; $col++) {
    $ws2.Cells.Item(241, $col).Value = 0
}
